# Updated cryptos list with GitHub Actions: refresh Price (D) and
# Volume(1h) (E) columns, and fix the row-41/row-42 coin data
# (Filecoin <-> Maker swapped back to correct rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that hold plain-decimal-looking text (e.g. "543.55") need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# and we lose the exact text representation (trailing zeros, etc).
$ws.Range("D2").Value = "61.973.33"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "2.993.97"
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.55"
$ws.Range("E5").Value = "  -3.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.95"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "3.008.08"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("E11").Value = "  -6.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.368"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "3.520.47"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "61.994.08"
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "3.002.03"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("E18").Value = "  -3.37%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.32"
$ws.Range("E21").Value = "  -5.69%  "
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.18"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").Value = "3.118.71"
$ws.Range("E25").Value = "  -4.30%  "
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "0.0₃0935"
$ws.Range("E29").Value = "  -6.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.29"
$ws.Range("E30").Value = "  -7.38%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.48"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "160.86"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.55"
$ws.Range("E39").Value = "  -4.94%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.421.24"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("E43").Value = "  -5.40%  "
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0592"
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.34"
$ws.Range("E46").Value = "  +4.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.996"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "269.35"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.60"
$ws.Range("E51").Value = "  -5.10%  "
